$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets("ALC")
$ws.Range("H17").Value = 1989.1346
$ws.Range("J17").Value = 1989.1346
$ws.Range("L17").Value = 5967.4038
$ws.Range("N17").Value = -6303.4038

$ws = $wb.Worksheets("ARM")
$ws.Range("H5").Value = 130.625
$ws.Range("I5").Value = 109.6
$ws.Range("K5").Value = 109.6
$ws.Range("M5").Value = 2.400000000000006

$ws.Range("H24").Value = 10000
$ws.Range("I24").Value = 10000
$ws.Range("K24").Value = 10000
$ws.Range("M24").Value = -9626

$ws.Range("H74").Value = 1520.3549
$ws.Range("I74").Value = 1620.0416
$ws.Range("K74").Value = 1620.0416
$ws.Range("M74").Value = -746.0416

$ws.Range("H77").Value = 1520.3549
$ws.Range("I77").Value = 1620.0416
$ws.Range("K77").Value = 8100.208000000001
$ws.Range("M77").Value = -3732.208000000001

$ws.Range("H100").Value = 10000
$ws.Range("I100").Value = 10000
$ws.Range("K100").Value = 10000
$ws.Range("M100").Value = -8918

$ws.Range("H132").Value = 3498.3125
$ws.Range("I132").Value = 3618.6897
$ws.Range("K132").Value = 10856.0691
$ws.Range("M132").Value = -8326.069100000001

$ws = $wb.Worksheets("BSM")
$ws.Range("H4").Value = 130.625
$ws.Range("I4").Value = 109.6
$ws.Range("K4").Value = 109.6
$ws.Range("M4").Value = 5.400000000000006

$ws.Range("H59").Value = 97250
$ws.Range("J59").Value = 97250
$ws.Range("L59").Value = 97250
$ws.Range("N59").Value = -98944

$ws.Range("H134").Value = 52386.668
$ws.Range("I134").Value = 5464.222
$ws.Range("J134").Value = 333921.34
$ws.Range("K134").Value = 16392.666
$ws.Range("L134").Value = 1001764.02
$ws.Range("M134").Value = -13857.666
$ws.Range("N134").Value = -1006834.02

$ws.Range("H139").Value = 99000
$ws.Range("J139").Value = 99000
$ws.Range("L139").Value = 99000
$ws.Range("N139").Value = -109280

$ws = $wb.Worksheets("CRP")
$ws.Range("H31").Value = 24161.043
$ws.Range("I31").Value = 1795.0869
$ws.Range("J31").Value = 46527
$ws.Range("K31").Value = 1795.0869
$ws.Range("L31").Value = 46527
$ws.Range("M31").Value = -1500.0869
$ws.Range("N31").Value = -47117

$ws.Range("H34").Value = 24161.043
$ws.Range("I34").Value = 1795.0869
$ws.Range("J34").Value = 46527
$ws.Range("K34").Value = 1795.0869
$ws.Range("L34").Value = 46527
$ws.Range("M34").Value = -1593.0869
$ws.Range("N34").Value = -46931

$ws.Range("H139").Value = 95000
$ws.Range("J139").Value = 95000
$ws.Range("L139").Value = 95000
$ws.Range("N139").Value = -105280

$ws = $wb.Worksheets("CUL")
$ws.Range("H39").Value = 9392.23
$ws.Range("I39").Value = 1033.3334
$ws.Range("K39").Value = 3100.0002
$ws.Range("M39").Value = -2806.0002

$ws.Range("H56").Value = 6999.25
$ws.Range("I56").Value = 6999.25
$ws.Range("K56").Value = 6999.25
$ws.Range("M56").Value = -6469.25

$ws.Range("H87").Value = 25199.2
$ws.Range("I87").Value = 13498.5
$ws.Range("J87").Value = 28124.375
$ws.Range("K87").Value = 40495.5
$ws.Range("L87").Value = 84373.125
$ws.Range("M87").Value = -39247.5
$ws.Range("N87").Value = -86869.125

$ws.Range("H90").Value = 25199.2
$ws.Range("I90").Value = 13498.5
$ws.Range("J90").Value = 28124.375
$ws.Range("K90").Value = 121486.5
$ws.Range("L90").Value = 253119.375
$ws.Range("M90").Value = -115246.5
$ws.Range("N90").Value = -265599.375

$ws.Range("H131").Value = 2501.9846
$ws.Range("I131").Value = 1838.9
$ws.Range("J131").Value = 2622.5454
$ws.Range("K131").Value = 5516.700000000001
$ws.Range("L131").Value = 7867.6362
$ws.Range("M131").Value = -476.7000000000007
$ws.Range("N131").Value = -17947.6362

$ws = $wb.Worksheets("GSM")
$ws.Range("H80").Value = 1056860
$ws.Range("I80").Value = 773911.3
$ws.Range("J80").Value = 1669915.4
$ws.Range("K80").Value = 773911.3
$ws.Range("L80").Value = 1669915.4
$ws.Range("M80").Value = -772913.3
$ws.Range("N80").Value = -1671911.4

$ws.Range("H83").Value = 1056860
$ws.Range("I83").Value = 773911.3
$ws.Range("J83").Value = 1669915.4
$ws.Range("K83").Value = 3869556.5
$ws.Range("L83").Value = 8349577
$ws.Range("M83").Value = -3864564.5
$ws.Range("N83").Value = -8359561

$ws.Range("H105").Value = 40000
$ws.Range("J105").Value = 40000
$ws.Range("L105").Value = 40000
$ws.Range("N105").Value = -46988

$ws.Range("H113").Value = 423451.4
$ws.Range("I113").Value = 835079.9399999999
$ws.Range("J113").Value = 11822.917
$ws.Range("K113").Value = 835079.9399999999
$ws.Range("L113").Value = 11822.917
$ws.Range("M113").Value = -832909.9399999999
$ws.Range("N113").Value = -16162.917

$ws.Range("H122").Value = 3798.75
$ws.Range("I122").Value = 3747
$ws.Range("J122").Value = 3954
$ws.Range("K122").Value = 11241
$ws.Range("L122").Value = 11862
$ws.Range("M122").Value = -8791
$ws.Range("N122").Value = -16762

$ws.Range("H132").Value = 90912.664
$ws.Range("I132").Value = 10628
$ws.Range("K132").Value = 31884
$ws.Range("M132").Value = -29354

$ws = $wb.Worksheets("LTW")
$ws.Range("H7").Value = 4717.478
$ws.Range("I7").Value = 4437.2104
$ws.Range("J7").Value = 6048.75
$ws.Range("K7").Value = 4437.2104
$ws.Range("L7").Value = 6048.75
$ws.Range("M7").Value = -4325.2104
$ws.Range("N7").Value = -6272.75

$ws.Range("H46").Value = 2248.3333
$ws.Range("I46").Value = 2209.7778
$ws.Range("K46").Value = 2209.7778
$ws.Range("M46").Value = -2021.7778

$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").ClearContents()

$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").ClearContents()

$ws.Range("H126").Value = 4717.478
$ws.Range("I126").Value = 4437.2104
$ws.Range("J126").Value = 6048.75
$ws.Range("K126").Value = 13311.6312
$ws.Range("L126").Value = 18146.25
$ws.Range("M126").Value = -10841.6312
$ws.Range("N126").Value = -23086.25

$ws.Range("H138").Value = 124666.336
$ws.Range("J138").Value = 124666.336
$ws.Range("L138").Value = 124666.336
$ws.Range("N138").Value = -134946.336

$ws = $wb.Worksheets("WVR")
$ws.Range("H4").Value = 3000136
$ws.Range("I4").Value = 3750120
$ws.Range("J4").Value = 200
$ws.Range("K4").Value = 3750120
$ws.Range("L4").Value = 200
$ws.Range("M4").Value = -3750007
$ws.Range("N4").Value = -426

$ws.Range("H96").Value = 333851.34
$ws.Range("I96").Value = 333851.34
$ws.Range("K96").Value = 333851.34
$ws.Range("M96").Value = -332478.34

$ws.Range("H122").Value = 29414630
$ws.Range("I122").Value = 43480388
$ws.Range("K122").Value = 130441164
$ws.Range("M122").Value = -130438714

$ws.Range("H139").Value = 100000
$ws.Range("J139").Value = 100000
$ws.Range("L139").Value = 100000
$ws.Range("N139").Value = -110280
